$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.559.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -3.31%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.606.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -2.01%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''571.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -4.45%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''154.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.84%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.04%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -3.18%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.604.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -2.00%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -7.82%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -0.62%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -5.01%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D14").Value = '''27.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -4.17%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.075.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -2.04%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '''  -7.77%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''63.486.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -3.25%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.595.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -1.82%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E20").Value = '''  +0.25%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  -6.33%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''340.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -3.83%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +0.04%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''67.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -3.89%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +1.26%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -6.15%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''9.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -6.68%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''577.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +2.84%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -4.46%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -0.04%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -2.06%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -3.73%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -4.53%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -5.77%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -1.73%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -3.16%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -5.31%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -0.03%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''19.61'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -4.68%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''153.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.16%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -5.82%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -0.03%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -3.21%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  -0.62%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''156.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.83%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''23.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.36%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -5.99%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -6.23%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -2.65%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.0998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.92%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -4.90%  '
$ws.Range("E51").Style = "Normal"
